# 12072110 added review dp
# Replace the "Sault" tracklist (Sheet1 / Sheet3) with the Dinner Party
# "Freeze Tag" / 9th Wonder remix tracklist, shrink the table from 15 to 7
# songs, update the header to a combined "Title/Composer" column, rename the
# web-query defined name / connections / query tables from "sault1" to
# "martinglasperwonderwashington1", and resize the columns to fit the much
# longer composer/performer text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New tracklist data (applies identically to Sheet1 and Sheet3, which are a
# duplicated web-query result table).
# ---------------------------------------------------------------------------
$titles = @(
    "Sleepless Nights",
    "Love You Bad",
    "From My Heart and My Soul",
    "First Responders",
    "The Mighty Tree",
    "Freeze Tag",
    "Luv U"
)

$composers = @(
    "Buddy / Patrick Douthit / Robert Glasper / Terrace Martin / Michael E. Neil / Reuben Vincent / Kamasi Washington",
    "Patrick Douthit / Robert Glasper / Terrace Martin / Michael E. Neil / Kamasi Washington / Malaya Watson",
    "Tarriona 'Tank' Ball / Patrick Douthit / Robert Glasper / Terrace Martin / Michael E. Neil / Kamasi Washington",
    "Patrick Douthit / Robert Glasper / Terrence Henderson / Terrace Martin / Bilal Oliver / Kamasi Washington",
    "Patrick Douthit / Marlanna Evans / Robert Glasper / Herbie Hancock / Terrace Martin / Kamasi Washington",
    "Patrick Douthit / Cordae Dunston / Robert Glasper / Terrace Martin / Michael E. Neil / Kamasi Washington",
    "Calvin Broadus / Patrick Douthit / Robert Glasper / Alex Isley / Terrace Martin / Kamasi Washington"
)

$performers = @(
    "Robert Glasper / Terrace Martin / Kamasi Washington / 9th Wonder feat. Buddy, Phoelix, Snoop Dogg",
    "Robert Glasper / Terrace Martin / Kamasi Washington / 9th Wonder feat. Malaya, Phoelix",
    "Robert Glasper / Terrace Martin / Kamasi Washington / 9th Wonder feat. Phoelix, Tank and the Bangas",
    "Robert Glasper / Terrace Martin / Kamasi Washington / 9th Wonder feat. Bilal, Punch",
    "Robert Glasper / Terrace Martin / Kamasi Washington / 9th Wonder feat. Rapsody, Herbie Hancock, Snoop Dogg",
    "Robert Glasper / Terrace Martin / Kamasi Washington / 9th Wonder feat. Phoelix, Cordae, Snoop Dogg",
    "Robert Glasper / Terrace Martin / Kamasi Washington / 9th Wonder feat. Snoop Dogg, Alex Isley"
)

$times = @(
    0.16597222222222222,
    0.09375,
    0.10555555555555556,
    0.12638888888888888,
    0.094444444444444442,
    0.13263888888888889,
    0.13749999999999998
)

foreach ($sheetName in @("Sheet1", "Sheet3")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Header row: drop the "No." / "Composer" headers, combine into
    # "Title/Composer" over B1, keep Performer/Time.
    $ws.Range("A1").Value = ""
    $ws.Range("B1").Value = "Title/Composer"
    $ws.Range("C1").Value = ""
    $ws.Range("D1").Value = "Performer"
    $ws.Range("E1").Value = "Time"

    for ($i = 0; $i -lt 7; $i++) {
        $r = $i + 2
        $ws.Range("A$r").Value = $i + 1
        $ws.Range("B$r").Value = $titles[$i]
        $ws.Range("C$r").Value = $composers[$i]
        $ws.Range("D$r").Value = $performers[$i]
        $ws.Range("E$r").Value = $times[$i]
    }

    # The old tracklist had 15 rows (2..16); clear the now-unused tail rows
    # (8 rows used, 9..16 blank) but keep their per-cell style/formatting.
    for ($r = 9; $r -le 16; $r++) {
        $ws.Range("A$r").Value = ""
        $ws.Range("B$r").Value = ""
        $ws.Range("C$r").Value = ""
        $ws.Range("D$r").Value = ""
        $ws.Range("E$r").Value = ""
    }

    # Column widths: A/E shrink, B widens slightly, C & D become equally
    # wide (one consistent width for both long-text columns).
    $ws.Columns.Item(1).ColumnWidth = 1.8571428571428572
    $ws.Columns.Item(2).ColumnWidth = 25.714285714285715
    $ws.Columns.Item(3).ColumnWidth = 80.14285714285714
    $ws.Columns.Item(4).ColumnWidth = 80.14285714285714
    $ws.Columns.Item(5).ColumnWidth = 5.0
}

# ---------------------------------------------------------------------------
# Rename the web-query defined name "sault1" -> "martinglasperwonderwashington1"
# and shrink its range from row 16 to row 12 (7 data rows + header).
# ---------------------------------------------------------------------------
foreach ($n in @($wb.Names)) {
    if ($n.Name -eq "Sheet1!sault1") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$E`$12"
        $n.Name = "martinglasperwonderwashington1"
    } elseif ($n.Name -eq "Sheet3!sault1") {
        $n.RefersTo = "=Sheet3!`$A`$1:`$E`$12"
        $n.Name = "martinglasperwonderwashington1"
    }
}

# ---------------------------------------------------------------------------
# Update the web-query connections' source URL (same rename).
# ---------------------------------------------------------------------------
foreach ($conn in $wb.Connections) {
    $conn.ODBCConnection
}

$newUrl = "http://bm.planetky.com/martinglasperwonderwashington1.htm"
foreach ($conn in @($wb.Connections)) {
    try { $conn.WorkbookConnection.Refresh } catch {}
}

# ---------------------------------------------------------------------------
# Update the two queryTables' names (drives queryTable1.xml / queryTable2.xml).
# ---------------------------------------------------------------------------
foreach ($sheetName in @("Sheet1", "Sheet3")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($qt in @($ws.QueryTables)) {
        $qt.Name = "martinglasperwonderwashington1"
        $qt.Connection = $newUrl
    }
}
